# Ok Alergenos en HTML
# Rebuild the data rows (2..12) of the "Tienda - Velázquez" closing report
# with the refreshed figures for 14/02/2025.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns H (total_arqueo_ciego) and J (total_operaciones) hold numeric-looking
# text (e.g. "488.10", "7") that must be stored as TEXT, not as numbers.
# Force the cells to text format before writing, then restore the default
# style so the cells don't end up flagged with an explicit style index.
$textColH = $ws.Range("H2:H12")
$textColJ = $ws.Range("J2:J12")
$textColH.NumberFormat = "@"
$textColJ.NumberFormat = "@"

$data = @(
    @("Tienda - Velázquez","V2","BAR","14/02/2025",8877,"Mañana","EUROS","488.10","38,00","7"),
    @("Tienda - Velázquez","V2","BAR","14/02/2025",8877,"Mañana","TARJETA VISA","597.95","100,55","20"),
    @("Tienda - Velázquez","V1","SERVIDOR TIENDA","14/02/2025",8876,"Mañana","EUROS","908.11","488,06","60"),
    @("Tienda - Velázquez","V1","SERVIDOR TIENDA","14/02/2025",8876,"Mañana","SMS","0.00","5,40","1"),
    @("Tienda - Velázquez","V1","SERVIDOR TIENDA","14/02/2025",8876,"Mañana","TARJETA VISA","1063.98","1084,13","105"),
    @("Tienda - Velázquez","V1","SERVIDOR TIENDA","14/02/2025",8879,"Mañana","EUROS","1199.50","830,60","82"),
    @("Tienda - Velázquez","V1","SERVIDOR TIENDA","14/02/2025",8879,"Mañana","TARJETA VISA","2825.69","2791,89","229"),
    @("Tienda - Velázquez","V2","BAR","14/02/2025",8878,"Mañana","EUROS","646.20","189,00","22"),
    @("Tienda - Velázquez","V2","BAR","14/02/2025",8878,"Mañana","TARJETA VISA","1226.39","468,34","49"),
    @("Tienda - Velázquez","V1","SERVIDOR TIENDA","14/02/2025",8880,"Mañana","EUROS","2161.42","685,51","66"),
    @("Tienda - Velázquez","V1","SERVIDOR TIENDA","14/02/2025",8880,"Mañana","TARJETA VISA","3045.82","1528,91","105")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
    $ws.Cells.Item($row, 8).Value = $r[7]
    $ws.Cells.Item($row, 9).Value = $r[8]
    $ws.Cells.Item($row, 10).Value = $r[9]
    $row = $row + 1
}

# Restore the default (unstyled) look for the cells we forced to text format.
$textColH.Style = "Normal"
$textColJ.Style = "Normal"
